$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 11.190867
$ws.Range("H2").Value = 33.572601
$ws.Range("I2").Value = 0.1514016037116739
$ws.Range("J2").Value = 0.1514016037116739
$ws.Range("M2").Value = 31.618405
$ws.Range("N2").Value = 94.855215
$ws.Range("O2").Value = 0.8578613706944929
$ws.Range("P2").Value = 0.8578613706944929
$ws.Range("Q2").Value = 353.837365107135
$ws.Range("R2").Value = 3184.536285964215
$ws.Range("S2").Value = 0.1298815872854409
$ws.Range("T2").Value = 0.1298815872854409

# Row 3
$ws.Range("G3").Value = 11.190867
$ws.Range("H3").Value = 33.572601
$ws.Range("I3").Value = 0.1514016037116739
$ws.Range("J3").Value = 0.1514016037116739
$ws.Range("O3").Value = 0.08747555172986397
$ws.Range("P3").Value = 0.08747555172986396
$ws.Range("Q3").Value = 36.080560091347
$ws.Range("R3").Value = 324.7250408221229
$ws.Range("S3").Value = 0.01324393881746489
$ws.Range("T3").Value = 0.01324393881746489

# Row 4
$ws.Range("G4").Value = 11.190867
$ws.Range("H4").Value = 33.572601
$ws.Range("I4").Value = 0.1514016037116739
$ws.Range("J4").Value = 0.1514016037116739
$ws.Range("M4").Value = 2.014730333333334
$ws.Range("N4").Value = 6.044191000000001
$ws.Range("O4").Value = 0.05466307757564324
$ws.Range("P4").Value = 0.05466307757564324
$ws.Range("Q4").Value = 22.546579201199
$ws.Range("R4").Value = 202.919212810791
$ws.Range("S4").Value = 0.008276077608768024
$ws.Range("T4").Value = 0.008276077608768024

# Row 5
$ws.Range("I5").Value = 0.2043341870182926
$ws.Range("J5").Value = 0.2043341870182926
$ws.Range("M5").Value = 31.618405
$ws.Range("N5").Value = 94.855215
$ws.Range("O5").Value = 0.8578613706944929
$ws.Range("P5").Value = 0.8578613706944929
$ws.Range("Q5").Value = 477.5449438009249
$ws.Range("R5").Value = 4297.904494208325
$ws.Range("S5").Value = 0.1752904057552574
$ws.Range("T5").Value = 0.1752904057552573

# Row 6
$ws.Range("I6").Value = 0.2043341870182926
$ws.Range("J6").Value = 0.2043341870182926
$ws.Range("O6").Value = 0.08747555172986397
$ws.Range("P6").Value = 0.08747555172986396
$ws.Range("S6").Value = 0.01787424574669836
$ws.Range("T6").Value = 0.01787424574669835

# Row 7
$ws.Range("I7").Value = 0.2043341870182926
$ws.Range("J7").Value = 0.2043341870182926
$ws.Range("M7").Value = 2.014730333333334
$ws.Range("N7").Value = 6.044191000000001
$ws.Range("O7").Value = 0.05466307757564324
$ws.Range("P7").Value = 0.05466307757564324
$ws.Range("S7").Value = 0.01116953551633692
$ws.Range("T7").Value = 0.01116953551633692

# Row 8
$ws.Range("G8").Value = 47.62086333333333
$ws.Range("I8").Value = 0.6442642092700336
$ws.Range("J8").Value = 0.6442642092700336
$ws.Range("M8").Value = 31.618405
$ws.Range("N8").Value = 94.855215
$ws.Range("O8").Value = 0.8578613706944929
$ws.Range("P8").Value = 0.8578613706944929
$ws.Range("Q8").Value = 1505.695743322983
$ws.Range("R8").Value = 13551.26168990685
$ws.Range("S8").Value = 0.5526893776537946
$ws.Range("T8").Value = 0.5526893776537946

# Row 9
$ws.Range("G9").Value = 47.62086333333333
$ws.Range("I9").Value = 0.6442642092700336
$ws.Range("J9").Value = 0.6442642092700336
$ws.Range("O9").Value = 0.08747555172986397
$ws.Range("P9").Value = 0.08747555172986396
$ws.Range("Q9").Value = 153.5347905662855
$ws.Range("S9").Value = 0.05635736716570073
$ws.Range("T9").Value = 0.05635736716570072

# Row 10
$ws.Range("G10").Value = 47.62086333333333
$ws.Range("I10").Value = 0.6442642092700336
$ws.Range("J10").Value = 0.6442642092700336
$ws.Range("M10").Value = 2.014730333333334
$ws.Range("N10").Value = 6.044191000000001
$ws.Range("O10").Value = 0.05466307757564324
$ws.Range("P10").Value = 0.05466307757564324
$ws.Range("Q10").Value = 95.94319785718777
$ws.Range("R10").Value = 863.4887807146899
$ws.Range("S10").Value = 0.0352174644505383
$ws.Range("T10").Value = 0.0352174644505383
